$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.393.01"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.607.26"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.40"
$ws.Range("E5").Value = "  +3.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.40"
$ws.Range("E6").Value = "  +0.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.57"
$ws.Range("E9").Value = "  -1.68%  "

$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("E11").Value = "  -2.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.371"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.066.71"
$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.46"
$ws.Range("E14").Value = "  +4.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.395.78"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("E16").Value = "  +1.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.610.74"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("E18").Value = "  +4.81%  "

$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.63"
$ws.Range("E20").Value = "  +0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.533"
$ws.Range("E23").Value = "  +3.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.66"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("E27").Value = "  +3.84%  "

$ws.Range("E28").Value = "  +11.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0800"
$ws.Range("E29").Value = "  +1.11%  "

$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.02"
$ws.Range("E32").Value = "  +4.55%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.33"
$ws.Range("E34").Value = "  +11.33%  "

$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.987"
$ws.Range("E36").Value = "  +2.59%  "

$ws.Range("E37").Value = "  +3.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.17"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("E39").Value = "  +3.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "311.77"
$ws.Range("E40").Value = "  +4.15%  "

$ws.Range("E41").Value = "  -1.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.96"
$ws.Range("E42").Value = "  -3.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0995"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.94"
$ws.Range("E45").Value = "  +2.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.608"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0552"
$ws.Range("E47").Value = "  +1.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.98"
$ws.Range("E48").Value = "  +4.12%  "

$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.92"
$ws.Range("E50").Value = "  +2.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.73"
$ws.Range("E51").Value = "  +0.46%  "
